# Applies the cryptos-list refresh described by the commit diff.
# D-column cells whose new value parses as a plain number get NumberFormat="@"
# forced (then reset to the "Normal" style) so Excel keeps storing them as the
# literal text the sheet already used (e.g. "4.30", "217.26") instead of silently
# coercing them into numeric cells and losing the formatting (trailing zeros, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.244.51'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +1.81%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.646.17'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  +0.56%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.18%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '217.26'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +0.83%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '0.507'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +0.54%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  -0.20%  '; ForceText = $false },
    @{ Cell = 'E8'; Value = '  -0.27%  '; ForceText = $false },
    @{ Cell = 'E9'; Value = '  -0.20%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '20.04'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +1.63%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.0793'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -0.33%  '; ForceText = $false },
    @{ Cell = 'B12'; Value = 'Polkadot'; ForceText = $false },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; ForceText = $false },
    @{ Cell = 'D12'; Value = '4.30'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  +0.44%  '; ForceText = $false },
    @{ Cell = 'B13'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false },
    @{ Cell = 'D13'; Value = '1.874.94'; ForceText = $false },
    @{ Cell = 'E13'; Value = '  +0.69%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '1.646.31'; ForceText = $false },
    @{ Cell = 'E14'; Value = '  +0.54%  '; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -1.40%  '; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -0.15%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '63.67'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  +0.89%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '26.238.79'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  +1.68%  '; ForceText = $false },
    @{ Cell = 'E19'; Value = '  -0.25%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '195.87'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +1.37%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '4.44'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -0.11%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  +0.77%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '6.36'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  +0.31%  '; ForceText = $false },
    @{ Cell = 'B24'; Value = 'Monero'; ForceText = $false },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false },
    @{ Cell = 'D24'; Value = '143.38'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  +0.49%  '; ForceText = $false },
    @{ Cell = 'B25'; Value = 'Toncoin'; ForceText = $false },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false },
    @{ Cell = 'D25'; Value = '1.78'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -2.17%  '; ForceText = $false },
    @{ Cell = 'E26'; Value = '  -0.15%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +2.14%  '; ForceText = $false },
    @{ Cell = 'E28'; Value = '  -0.10%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '15.62'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  +0.40%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  +1.31%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '0.0503'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  +1.54%  '; ForceText = $false },
    @{ Cell = 'E32'; Value = '  +0.59%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  +0.45%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  +1.53%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  +1.19%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '0.915'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  +1.43%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '1.139.38'; ForceText = $false },
    @{ Cell = 'E37'; Value = '  +0.60%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.556'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +2.22%  '; ForceText = $false },
    @{ Cell = 'E39'; Value = '  -0.70%  '; ForceText = $false },
    @{ Cell = 'E40'; Value = '  +0.60%  '; ForceText = $false },
    @{ Cell = 'E41'; Value = '  -0.10%  '; ForceText = $false },
    @{ Cell = 'E42'; Value = '  +1.84%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '100.46'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -0.24%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '0.799'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -1.04%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '1.784.40'; ForceText = $false },
    @{ Cell = 'E45'; Value = '  +0.71%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '56.37'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +1.97%  '; ForceText = $false },
    @{ Cell = 'B47'; Value = 'RenderToken'; ForceText = $false },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false },
    @{ Cell = 'D47'; Value = '1.48'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +5.00%  '; ForceText = $false },
    @{ Cell = 'B48'; Value = 'Cronos'; ForceText = $false },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; ForceText = $false },
    @{ Cell = 'D48'; Value = '0.0518'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +2.92%  '; ForceText = $false },
    @{ Cell = 'B49'; Value = 'EnergySwap'; ForceText = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false },
    @{ Cell = 'D49'; Value = '7.73'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  +2.98%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'Mantle'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.418'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.0978'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +2.97%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Type a number-looking literal as text: set a text format first so the
        # assignment doesn't get reinterpreted as a numeric value, then drop the
        # explicit format again (style resets to the default "Normal") so the
        # cell is indistinguishable from one that was always plain text.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}

